{"js": "// Office.js (Word JavaScript API) edit script.\n// Body of: async (context) => { ... }\n//\n// Target change (per the diff):\n//   1. \"## [1] 97.46328\"  ->  \"## [1] 97.39653\"\n//   2. \"##      0 747  18\" ->  \"##      0 746  18\"\n//   3. \"##      1  20 713\" ->  \"##      1  21 713\"\n\nconst body = context.document.body;\n\nasync function replaceUnique(oldText, newText) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length !== 1) {\n    throw new Error(\n      `Expected exactly 1 match for \"${oldText}\", found ${results.items.length}`\n    );\n  }\n\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 1) Logistic-regression accuracy output.\nawait replaceUnique(\"97.46328\", \"97.39653\");\n\n// 2) Confusion-matrix counts for the knn.38 table.\nawait replaceUnique(\"0 747  18\", \"0 746  18\");\nawait replaceUnique(\"1  20 713\", \"1  21 713\");\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# $word.ActiveDocument is the open document.\n#\n# Target change (per the diff):\n#   1. \"## [1] 97.46328\"   -> \"## [1] 97.39653\"\n#   2. \"##      0 747  18\" -> \"##      0 746  18\"\n#   3. \"##      1  20 713\" -> \"##      1  21 713\"\n\n$d = $word.ActiveDocument\n\n# Find.Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards,\n#              MatchSoundsLike, MatchAllWordForms, Forward, Wrap,\n#              Format, ReplaceWith, Replace)\n#   Wrap    = 1 (wdFindContinue)\n#   Replace = 1 (wdReplaceOne) -> each search string is unique in the\n#              document (checked below), so replacing the first hit is\n#              sufficient and avoids touching any unrelated text.\n\nfunction Count-Occurrences($haystack, $needle) {\n    return ($haystack.Split($needle)).Length - 1\n}\n\nfunction Replace-UniqueText($doc, $oldText, $newText) {\n    $before = Count-Occurrences $doc.Content.Text $oldText\n    if ($before -ne 1) {\n        throw \"Expected exactly 1 occurrence of '$oldText', found $before\"\n    }\n\n    $find = $doc.Content.Find\n    $find.ClearFormatting()\n    $ok = $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 1)\n    if (-not $ok) {\n        throw \"Find/Replace failed for '$oldText'\"\n    }\n\n    $after = Count-Occurrences $doc.Content.Text $newText\n    if ($after -lt 1) {\n        throw \"Replacement text '$newText' not found after replace\"\n    }\n}\n\n# 1) Logistic-regression accuracy output.\nReplace-UniqueText $d \"97.46328\" \"97.39653\"\n\n# 2) Confusion-matrix counts for the knn.38 table.\nReplace-UniqueText $d \"0 747  18\" \"0 746  18\"\nReplace-UniqueText $d \"1  20 713\" \"1  21 713\"\n"}
